$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The worksheet is protected; unprotect so values can be written, then
# restore protection afterwards to preserve the original workbook state.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A42).
$confidentialText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."
$ws.Range("A42").Value = $confidentialText

# Update the Weight (D) and Percent Change (E) values for each holding row.
$ws.Range("D2").Value = 0.05803146228699353
$ws.Range("E2").Value = -0.02493844809784773
$ws.Range("D3").Value = 0.05278442504678097
$ws.Range("E3").Value = -0.02936279088656935
$ws.Range("D4").Value = 0.3076215324782955
$ws.Range("E4").Value = -0.007317073170731714
$ws.Range("D5").Value = 0.03455554476557846
$ws.Range("E5").Value = -0.02232382417623313
$ws.Range("D6").Value = 0.03141855543072052
$ws.Range("E6").Value = -0.02150966826326051
$ws.Range("D7").Value = 0.03058765573742026
$ws.Range("E7").Value = -0.006875236533366968
$ws.Range("D8").Value = 0.02896229931982414
$ws.Range("E8").Value = -0.004026527711984862
$ws.Range("D9").Value = 0.02393231211559367
$ws.Range("E9").Value = -0.02586886420637768
$ws.Range("D10").Value = 0.02433168418180068
$ws.Range("E10").Value = -0.03075249112358258
$ws.Range("D11").Value = 0.02299883618306632
$ws.Range("E11").Value = -0.01298404723844326
$ws.Range("D12").Value = 0.02316958221137222
$ws.Range("E12").Value = -0.009381765696415578
$ws.Range("D13").Value = 0.02141956940274047
$ws.Range("E13").Value = -0.02102713711675019
$ws.Range("D14").Value = 0.02202495142902528
$ws.Range("E14").Value = -0.01058471421271645
$ws.Range("D15").Value = 0.0212651155252836
$ws.Range("E15").Value = -0.04131131676041855
$ws.Range("D16").Value = 0.02206407402810353
$ws.Range("E16").Value = -0.01537527325722621
$ws.Range("D17").Value = 0.01970117622898018
$ws.Range("E17").Value = -0.02989581349800075
$ws.Range("D18").Value = 0.01415873656942611
$ws.Range("E18").Value = -0.02579941860465129
$ws.Range("D19").Value = 0.0172568176756117
$ws.Range("E19").Value = -0.01530434782608703
$ws.Range("D20").Value = 0.01568655445233399
$ws.Range("E20").Value = -0.002220703792278789
$ws.Range("D21").Value = 0.01688531376217275
$ws.Range("E21").Value = -0.009077405512460945
$ws.Range("D22").Value = 0.01323094145265533
$ws.Range("E22").Value = -0.04424821775761512
$ws.Range("D23").Value = 0.01513798058361485
$ws.Range("E23").Value = -0.005154639175257714
$ws.Range("D24").Value = 0.01486851698338822
$ws.Range("E24").Value = -0.008989460632361945
$ws.Range("D25").Value = 0.01415648568016407
$ws.Range("E25").Value = -0.02328222600795005
$ws.Range("D26").Value = 0.01395862179551078
$ws.Range("E26").Value = -0.01723118506630628
$ws.Range("D27").Value = 0.0132275115261608
$ws.Range("E27").Value = -0.02720245040840152
$ws.Range("D28").Value = 0.01372399338624424
$ws.Range("E28").Value = 0.006279287722586568
$ws.Range("D29").Value = 0.01434030830323037
$ws.Range("E29").Value = 0.008640406607369755
$ws.Range("D30").Value = 0.01332354946800768
$ws.Range("E30").Value = -0.01209936928819666
$ws.Range("D31").Value = 0.01248921984821289
$ws.Range("E31").Value = -0.009749399244764922
$ws.Range("D32").Value = 0.01354595876413746
$ws.Range("E32").Value = 0.0006092784402471629
$ws.Range("D33").Value = 0.01245384873123803
$ws.Range("E33").Value = 0.006971340046475705
$ws.Range("D34").Value = 0.006133673239048941
$ws.Range("E34").Value = -0.03828746177370024
$ws.Range("D35").Value = 0.005306525027852074
$ws.Range("E35").Value = -0.02040074331421182
$ws.Range("D36").Value = 0.005329248290878346
$ws.Range("E36").Value = -0.03495575221238933
$ws.Range("D37").Value = 0.00520051886213046
$ws.Range("E37").Value = -0.02699973206372763
$ws.Range("D38").Value = 0.004716899226401516
$ws.Range("E38").Value = -0.04019815029427132
$ws.Range("D39").Value = 0.9999999999999999
$ws.Range("E39").Value = -0.01477751674608019

# Restore sheet protection to match the original workbook.
$ws.Protect("D382")
